$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title cell value
$ws.Range("A1").Value = "CreateContact"

# Move the active selection to A2 (matches the commit's sheetView change)
$ws.Range("A2").Select()
